# Update countries & provincias Spain
#
# Refresh of the COVID-19 "Pais" dashboard data: a handful of countries'
# totals were updated, which in a couple of cases causes them to swap
# rank/row order with their neighbour, plus the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "last updated" footer -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Septiembre de 2020 a las 21:23"

# --- Estados Unidos (row 4) --------------------------------------------------
$ws.Cells.Item(4, 2).Value = 6853787
$ws.Cells.Item(4, 3).Value = 25486
$ws.Cells.Item(4, 4).Value = 4137645
$ws.Cells.Item(4, 5).Value = 2514305
$ws.Cells.Item(4, 7).Value = 489
$ws.Cells.Item(4, 8).Value = 201837

# --- India (row 5) -----------------------------------------------------------
$ws.Cells.Item(5, 2).Value = 5212623
$ws.Cells.Item(5, 3).Value = 96730
$ws.Cells.Item(5, 4).Value = 4107515
$ws.Cells.Item(5, 5).Value = 1020708
$ws.Cells.Item(5, 7).Value = 1170
$ws.Cells.Item(5, 8).Value = 84400

# --- Alemania (row 25) --------------------------------------------------------
$ws.Cells.Item(25, 2).Value = 268258
$ws.Cells.Item(25, 3).Value = 1393
$ws.Cells.Item(25, 5).Value = 19706

# --- Costa Rica overtakes Nepal (rows 55-56 swap) -----------------------------
$ws.Cells.Item(55, 1).Value = "Costa Rica"
$ws.Cells.Item(55, 2).Value = 60818
$ws.Cells.Item(55, 3).Value = 1302
$ws.Cells.Item(55, 4).Value = 22662
$ws.Cells.Item(55, 5).Value = 37490
$ws.Cells.Item(55, 7).Value = 17
$ws.Cells.Item(55, 8).Value = 666

$ws.Cells.Item(56, 1).Value = "Nepal"
$ws.Cells.Item(56, 2).Value = 59573
$ws.Cells.Item(56, 3).Value = 1246
$ws.Cells.Item(56, 4).Value = 42949
$ws.Cells.Item(56, 5).Value = 16241
$ws.Cells.Item(56, 7).Value = 4
$ws.Cells.Item(56, 8).Value = 383

# --- Malaui (row 112) ----------------------------------------------------------
$ws.Cells.Item(112, 2).Value = 5711
$ws.Cells.Item(112, 3).Value = 7
$ws.Cells.Item(112, 4).Value = 4000
$ws.Cells.Item(112, 5).Value = 1532
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = 179

# --- Reunion overtakes Bahamas (rows 140-141 swap) ------------------------------
$ws.Cells.Item(140, 1).Value = "Reunion"
$ws.Cells.Item(140, 2).Value = 3099
$ws.Cells.Item(140, 3).Value = 97
$ws.Cells.Item(140, 4).Value = 1794
$ws.Cells.Item(140, 5).Value = 1290
$ws.Cells.Item(140, 8).Value = 15

$ws.Cells.Item(141, 1).Value = "Bahamas"
$ws.Cells.Item(141, 2).Value = 3087
$ws.Cells.Item(141, 4).Value = 1533
$ws.Cells.Item(141, 5).Value = 1485
$ws.Cells.Item(141, 8).Value = 69

# --- Curazao overtakes Islas Caimanes (rows 188-189 swap) -----------------------
$ws.Cells.Item(188, 1).Value = "Curazao"
$ws.Cells.Item(188, 2).Value = 210
$ws.Cells.Item(188, 3).Value = 18
$ws.Cells.Item(188, 4).Value = 74
$ws.Cells.Item(188, 5).Value = 135

$ws.Cells.Item(189, 1).Value = "Islas Caimanes"
$ws.Cells.Item(189, 2).Value = 208
$ws.Cells.Item(189, 4).Value = 204
$ws.Cells.Item(189, 5).Value = 3

# --- Monaco overtakes Barbados (rows 190-191 swap) ------------------------------
$ws.Cells.Item(190, 1).Value = "Monaco"
$ws.Cells.Item(190, 2).Value = 186
$ws.Cells.Item(190, 3).Value = 5
$ws.Cells.Item(190, 4).Value = 147
$ws.Cells.Item(190, 5).Value = 38
$ws.Cells.Item(190, 8).Value = 1

$ws.Cells.Item(191, 1).Value = "Barbados"
$ws.Cells.Item(191, 2).Value = 185
$ws.Cells.Item(191, 4).Value = 171
$ws.Cells.Item(191, 5).Value = 7
$ws.Cells.Item(191, 8).Value = 7

# --- Islas Malvinas / Montserrat swap (rows 214-215) ----------------------------
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
